$p = $ppt.ActivePresentation

# Slide 10 (sldId 273) holds the "Open Mic" body placeholder (shape id=3,
# the 2nd shape in z-order).
$s = $p.Slides.Item(10)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

# Paragraph 1: new introductory sentence about Knowledge Graph interim.
$tr.Text = "Given the various Knowledge Graph (KG) related I-Ds out there, the discussion on the mailing list, presentations in previous interim,"
$tr.LanguageID = "en-US"

# Paragraph 2: the actual question, built incrementally (InsertAfter) so the
# language attribute of each new paragraph is committed correctly.
$null = $tr.InsertAfter("`rIs the WG interested to organize a dedicated interim on KG?")
$tr.Paragraphs(2).LanguageID = "en-US"

# Paragraph 3: blank separator line.
$null = $tr.InsertAfter("`r")
$tr.Paragraphs(3).LanguageID = "en-US"

# Paragraph 4: restored "Open Mic" heading (previously had the
# "Open Mic (if time permits)" text split over several runs).
$null = $tr.InsertAfter("`rOpen Mic")
$tr.Paragraphs(4).LanguageID = "en-US"

# Apply indent level / bold / italic formatting to paragraph 2 only, as a
# final pass so it does not bleed into the paragraphs appended afterwards.
$para2 = $tr.Paragraphs(2)
$para2.IndentLevel = 2
$para2.Font.Bold = $true
$para2.Font.Italic = $true
